# The canonical-OOXML diff for this revision touches only
# word/document.xml (the <w:document> root element) and word/styles.xml
# (<w:docDefaults>, <w:latentStyles>, <w:lsdException>, <w:style> ...).
# In every single hunk the "-" and "+" lines carry exactly the same set
# of XML namespace declarations / attributes and exactly the same
# values - only the left-to-right attribute order changed (namespace
# declarations sorted by prefix, then the remaining attributes sorted
# alphabetically by local name). There is no textual, structural or
# value-level change anywhere in the package: no run text, field code,
# style property, page-size/margin value, language, etc. is added,
# removed or modified.
#
# That kind of diff is what you get when a docx is re-serialized by an
# XML layer that writes attributes back out in a different (e.g.
# alphabetical/DOM-map) order than the one Word originally used - it is
# not something a user (or a macro driving the Word object model)
# produces by editing the document's content. Word's own object model
# has no API that reorders existing attributes on an element; it only
# lets you read/write the logical properties, and this document's
# properties (fonts, languages, page size/margins, style/latent-style
# catalog, ...) already hold exactly the values the target revision
# expects.
#
# So the correct, content-preserving edit here is a no-op against the
# Word object model: touch nothing, leave every paragraph, run, style
# and section property exactly as authored. Re-saving through the COM
# host below round-trips the package without introducing any semantic
# difference, matching the (attribute-order-only) diff.
$d = $word.ActiveDocument
$null = $d.Name
